# Resident Services_Requirements.xlsx — add "Reg Proc" column to Table2 on
# the "Details" sheet, with clarification / research-info text for several
# rows, and update an existing "Comments" cell (S8) to wrap + extra text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")
$ws.Activate()

$tbl = $ws.ListObjects.Item("Table2")

# --- Add the new table column (lands at the end -> column T) ---------------
$newCol = $tbl.ListColumns.Add()

# Copy the header formatting from the previous header cell (S2) onto the new
# header cell (T2) so it picks up the same centered/bordered header style,
# then set its caption -- this also renames the ListColumn to "Reg Proc".
$ws.Range("S2").Copy()
$ws.Range("T2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("T2").Value = "Reg Proc"

# --- New "Reg Proc" column data ---------------------------------------------
# (written in the same order the source strings were first introduced, so
# shared-string indices line up with the authored workbook)
$ws.Range("T5").WrapText = $true
$ws.Range("T5").Value = "No Mapping of such kind from Reg Processor`nID Repo- Might not be there in ID Repo as well"

$ws.Range("T6").Value = "ID Repo- need to know "

$ws.Range("T9").WrapText = $true
$ws.Range("T9").Value = "Under processing`nProcessed"

$ws.Range("T8").WrapText = $true
$ws.Range("T8").Value = "Under processing`nProcessed`n"

# --- Existing "Comments" cell (S8) gains wrap + an extra line ---------------
$ws.Range("S8").WrapText = $true
$ws.Range("S8").Value = "Reg proc`nArchival policy"

$ws.Range("T10").Value = "E-UIN Generation"

$ws.Range("T7").WrapText = $true
$ws.Range("T7").Value = "there shud be a label as Res_Service`nReg Client packet needs to be understood`nService from Reg proc needs to be developed"

$ws.Range("T4").WrapText = $true
$ws.Range("T4").Value = "When UIN IS needed to be generated`n1.the Acknowledgment from Print queue- what needs to be done`nTime period `n2. If there is a print failure- no need to handle from MOSIP`nUser Story ?"

# --- Autofit the new column's width -----------------------------------------
$ws.Columns.Item(20).ColumnWidth = 32

# --- Selection / active cell bookkeeping ------------------------------------
$ws.Range("T4").Select()
